# Insert a new weekly price-observation row for "Papa" (Patagonia, 1a (cosecha))
# at row 606 of the "Feria Lagunitas de Puerto Montt" sheet. All existing rows
# from 606 onward shift down by one (old row 645 becomes row 646), and the
# sheet's used range grows from A1:R645 to A1:R646.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 606..645 down to 607..646, leaving a fresh blank row 606.
$ws.Rows.Item(606).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A606").Value2 = 4
$ws.Range("B606").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C606").Value2 = "Los Lagos"
$ws.Range("D606").Value2 = 45021
$ws.Range("E606").Value2 = 10
$ws.Range("F606").Value2 = 100114001
$ws.Range("G606").Value2 = "Papa"
$ws.Range("H606").Value2 = "Patagonia"
$ws.Range("I606").Value2 = "1a (cosecha)"
$ws.Range("J606").Value2 = 150
$ws.Range("K606").Value2 = 12000
$ws.Range("L606").Value2 = 12000
$ws.Range("M606").Value2 = 12000
$ws.Range("N606").Value2 = "$/saco 25 kilos"
$ws.Range("O606").Value2 = "Provincia de Llanquihue"
$ws.Range("P606").Value2 = 480
$ws.Range("Q606").Value2 = 25
$ws.Range("R606").Value2 = "Hortaliza"
